$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "Datos actualizados" timestamp, 19:22 -> 19:52
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 19:52"

# Row 4: Estados Unidos - updated totals
$ws.Range("B4").Value = 521714
$ws.Range("C4").Value = 18838
$ws.Range("D4").Value = 28580
$ws.Range("E4").Value = 473070
$ws.Range("F4").Value = 10952
$ws.Range("G4").Value = 1317
$ws.Range("H4").Value = 20064

# Row 7: Francia - updated totals
$ws.Range("B7").Value = 129654
$ws.Range("C7").Value = 4785
$ws.Range("D7").Value = 26391
$ws.Range("E7").Value = 89431
$ws.Range("F7").Value = 6883
$ws.Range("G7").Value = 635
$ws.Range("H7").Value = 13832

# Row 8: Alemania - updated totals
$ws.Range("B8").Value = 123878
$ws.Range("C8").Value = 1707
$ws.Range("E8").Value = 67229

# Row 14: Suiza - updated totals
$ws.Range("B14").Value = 25107
$ws.Range("C14").Value = 556
$ws.Range("E14").Value = 12971

# Row 16: Canada - updated totals
$ws.Range("E16").Value = 15962
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = 600

# Row 24: India - updated totals
$ws.Range("D24").Value = 969
$ws.Range("E24").Value = 7082
$ws.Range("G24").Value = 39
$ws.Range("H24").Value = 288

# Peru moves up in the country list: it now sits right after Chile (row 27),
# before Noruega, so rows 28-34 each shift down to the next country and get
# that country's refreshed totals; Peru itself lands on row 28 with new data.

# Row 28: was Noruega -> now Peru (new data)
$ws.Range("A28").Value = "Peru"
$ws.Range("B28").Value = 6848
$ws.Range("C28").Value = 951
$ws.Range("D28").Value = 1739
$ws.Range("E28").Value = 4928
$ws.Range("F28").Value = 142
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 181

# Row 29: was Polonia -> now Noruega
$ws.Range("A29").Value = "Noruega"
$ws.Range("B29").Value = 6403
$ws.Range("C29").Value = 89
$ws.Range("D29").Value = 32
$ws.Range("E29").Value = 6254
$ws.Range("F29").Value = 67
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 117

# Row 30: was Australia -> now Polonia
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 6356
$ws.Range("C30").Value = 401
$ws.Range("D30").Value = 375
$ws.Range("E30").Value = 5773
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 27
$ws.Range("H30").Value = 208

# Row 31: was Japon -> now Australia
$ws.Range("A31").Value = "Australia"
$ws.Range("B31").Value = 6303
$ws.Range("C31").Value = 65
$ws.Range("D31").Value = 3265
$ws.Range("E31").Value = 2982
$ws.Range("F31").Value = 80
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 56

# Row 32: was Dinamarca -> now Japon
$ws.Range("A32").Value = "Japon"
$ws.Range("B32").Value = 6005
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 762
$ws.Range("E32").Value = 5144
$ws.Range("F32").Value = 109
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 99

# Row 33: was Rumania -> now Dinamarca
$ws.Range("A33").Value = "Dinamarca"
$ws.Range("B33").Value = 5996
$ws.Range("C33").Value = 177
$ws.Range("D33").Value = 1955
$ws.Range("E33").Value = 3781
$ws.Range("F33").Value = 106
$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 260

# Row 34: was Peru (old data) -> now Rumania (refreshed totals)
$ws.Range("A34").Value = "Rumania"
$ws.Range("B34").Value = 5990
$ws.Range("C34").Value = 523
$ws.Range("D34").Value = 758
$ws.Range("E34").Value = 4941
$ws.Range("F34").Value = 208
$ws.Range("G34").Value = 21
$ws.Range("H34").Value = 291

# Row 35 (Chequia) is unaffected by the Peru re-shuffle - no change needed.

# Row 36: Pakistan - updated totals
$ws.Range("B36").Value = 5011
$ws.Range("C36").Value = 316
$ws.Range("E36").Value = 4172

# Row 75: Kazajistan - updated totals
$ws.Range("B75").Value = 865
$ws.Range("C75").Value = 53
$ws.Range("E75").Value = 774

# Row 109: Estado de Palestina - updated totals
$ws.Range("D109").Value = 57
$ws.Range("E109").Value = 209
